# Apply updated task-order timestamps / filenames and rename sheets accordingly.
$wb = $excel.ActiveWorkbook

# --- Rename sheets (by index, matching original sheet order) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16512554921067002"
$wb.Worksheets.Item(2).Name = "NB_TO-16512554967732973"
$wb.Worksheets.Item(3).Name = "RS_TO-16512554967742703"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512554968212702"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512554968982954"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1651255492077702.csv"
$ws1.Range("B3").Value = "GNG_stims-16512554920897021.csv"
$ws1.Range("B4").Value = "go_stims-16512554920907385.csv"
$ws1.Range("B5").Value = "GNG_stims-16512554921047344.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16512554950632627.csv"
$ws2.Range("B3").Value = "ZB-match_7-16512554929511437.csv"
$ws2.Range("B4").Value = "OB-16512554936942945.csv"
$ws2.Range("B5").Value = "ZB-match_4-16512554934351454.csv"
$ws2.Range("B6").Value = "TB-16512554966512594.csv"
$ws2.Range("B7").Value = "TB-16512554957152936.csv"
$ws2.Range("B8").Value = "OB-16512554953052578.csv"
$ws2.Range("B9").Value = "TB-16512554967482963.csv"
$ws2.Range("B10").Value = "ZB-match_5-16512554928201795.csv"

# --- Sheet 3: RS_TO (no cell data changes, name already updated above) ---

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1651255496788298.csv"
$ws4.Range("B3").Value = "ZM_stims-16512554967772596.csv"
$ws4.Range("B4").Value = "MM_stims-16512554968042595.csv"
$ws4.Range("B5").Value = "ZM_stims-16512554967892687.csv"
$ws4.Range("B6").Value = "MM_stims-1651255496820294.csv"
$ws4.Range("B7").Value = "ZM_stims-16512554968042595.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1651255496825261.csv"
$ws5.Range("B3").Value = "vSAT_stims-16512554968832636.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512554968672976.csv"
$ws5.Range("B5").Value = "SAT_stims-16512554968513007.csv"

$wb.Save()
